$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update values on Tabelle1
$ws.Range("B2").Value = "nein"
$ws.Range("B3").Value = 7.3
$ws.Range("B4").Value = 7.3

# Move the active selection from B5 to B4
$ws.Range("B4").Select()
